$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers (e.g. "241.12").
# Excel would auto-convert these to numeric values if assigned directly,
# losing the original text formatting (trailing zeros, etc). Force them
# to remain text by temporarily setting a text number format.
$textForced = @{
    'D5' = '0.7918'
    'D6' = '241.12'
    'D8' = '0.3155'
    'D9' = '25.53'
    'D10' = '0.06996'
    'D11' = '0.08039'
    'D14' = '5.286'
    'D15' = '92.23'
    'D17' = '13.78'
    'D18' = '5.914'
    'D19' = '243.33'
    'D20' = '0.000007660'
    'D21' = '1.000'
    'D23' = '8.108'
    'D25' = '0.1646'
    'D26' = '9.278'
    'D27' = '164.17'
    'D28' = '18.60'
    'D30' = '1.393'
    'D31' = '1.531'
    'D32' = '4.374'
    'D33' = '0.05675'
    'D34' = '4.038'
    'D36' = '0.7330'
    'D37' = '0.9994'
    'D38' = '2.593'
    'D39' = '0.01899'
    'D40' = '2.775'
    'D41' = '0.4385'
    'D42' = '72.17'
    'D43' = '5.809'
    'D44' = '1.000'
    'D45' = '0.8373'
    'D46' = '102.34'
    'D48' = '1.855'
    'D49' = '9.855'
    'D50' = '7.424'
}
foreach ($addr in $textForced.Keys) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $textForced[$addr]
    $c.ClearFormats()
}

# Remaining cells already stay text when assigned normally (they contain
# characters, multiple dots, percent signs and padding spaces that keep
# Excel from reinterpreting them as numbers).
$plainValues = @{
    'D2' = '29.763.33'
    'E2' = '  -0.88%  '
    'D3' = '1.889.56'
    'E3' = '  -1.03%  '
    'E4' = '  +0.08%  '
    'E5' = '  -4.47%  '
    'E6' = '  -0.41%  '
    'E7' = '  +0.06%  '
    'E8' = '  -2.14%  '
    'E9' = '  -4.36%  '
    'E10' = '  -0.18%  '
    'E11' = '  +0.21%  '
    'E12' = '  +0.90%  '
    'D13' = '1.895.53'
    'E13' = '  -0.74%  '
    'E14' = '  +1.44%  '
    'E15' = '  -0.63%  '
    'D16' = '29.787.69'
    'E16' = '  -0.77%  '
    'E17' = '  -2.53%  '
    'E18' = '  +0.24%  '
    'E19' = '  -0.66%  '
    'E20' = '  -1.51%  '
    'E21' = '  +0.05%  '
    'D22' = '2.152.53'
    'E22' = '  -0.36%  '
    'E23' = '  +16.44%  '
    'E24' = '  +0.14%  '
    'E25' = '  +1.78%  '
    'E26' = '  +0.41%  '
    'E27' = '  -2.99%  '
    'E28' = '  -1.84%  '
    'E29' = '  -2.23%  '
    'E30' = '  +1.70%  '
    'E31' = '  +1.01%  '
    'E33' = '  +1.09%  '
    'E34' = '  -1.26%  '
    'E35' = '  -1.05%  '
    'E37' = '  +0.07%  '
    'E38' = '  -4.43%  '
    'E39' = '  -1.30%  '
    'E40' = '  -0.55%  '
    'E41' = '  -1.08%  '
    'E42' = '  -0.36%  '
    'E43' = '  -2.95%  '
    'E44' = '  +0.07%  '
    'E45' = '  -0.64%  '
    'E46' = '  +1.28%  '
    'D47' = '1.017.93'
    'E47' = '  +3.29%  '
    'E48' = '  -2.09%  '
    'E49' = '  +1.25%  '
    'E50' = '  -2.33%  '
    'D51' = '2.056.01'
    'E51' = '  -0.30%  '
}
foreach ($addr in $plainValues.Keys) {
    $ws.Range($addr).Value = $plainValues[$addr]
}
